# The workbook had two sheets ("futbol" and "futbol varones") both listing
# the same header row of person-related fields. Consolidate them into a
# single "Personas" sheet: keep/rename "futbol" -> "Personas" and remove the
# duplicate "futbol varones" sheet.

$wb = $excel.ActiveWorkbook

$sheetToKeep = $wb.Worksheets.Item("futbol")
$sheetToRemove = $wb.Worksheets.Item("futbol varones")

$excel.DisplayAlerts = $false
$sheetToRemove.Delete()
$excel.DisplayAlerts = $true

$sheetToKeep.Name = "Personas"
